# Merge the first two runs of the "YouTube / Quora / RosettaCode" line on
# slide 2 into a single run, dropping the leading tab character, to match
# the author's edit: "\tYouTube • Quora " + "• " -> "YouTube • Quora • "
# (keeping the second run's formatting/rPr), while leaving the trailing
# "RosettaCode" run untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(5)
$tr = $sh.TextFrame.TextRange

$bullet = [char]0x2022
$para = $tr.Paragraphs(2)

$run1 = $para.Runs(1, 1)
$run2 = $para.Runs(2, 1)

# Replace the second run's text with the fully merged text, keeping that
# run's original formatting (it already carries dirty="0").
$run2.Text = "YouTube " + $bullet + " Quora " + $bullet + " "

# Clear the first run's text so it collapses away entirely instead of
# leaving a stray empty run behind.
$para2 = $tr.Paragraphs(2)
$run1b = $para2.Runs(1, 1)
$run1b.Text = ""
